$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move row 12 ("CO_NW_11") to the bottom of the list (row 33), shifting rows 13-33 up by one ---

# 1) Stash row 12's original values (C/D/E - the columns that actually vary row to row)
$origC12 = $ws.Range("C12").Value()
$origD12 = $ws.Range("D12").Value()
$origE12 = $ws.Range("E12").Value()

# 2) Stash row 12's original formatting onto an unused scratch row so we can re-apply
#    it later to row 33 (PasteSpecial formats-only keeps the cellXfs in sync with the
#    data move, matching what Excel itself does on a drag/cut-insert of the row).
$ws.Range("A12:E12").Copy()
$ws.Range("A200:E200").PasteSpecial(-4122)

# 3) Shift rows 13..33 up into 12..32 - both formatting and the C/D/E values.
for ($n = 12; $n -le 32; $n++) {
    $src = $n + 1

    $ws.Range("A" + $src + ":E" + $src).Copy()
    $ws.Range("A" + $n + ":E" + $n).PasteSpecial(-4122)

    $cVal = $ws.Range("C" + $src).Value()
    $dVal = $ws.Range("D" + $src).Value()
    $eVal = $ws.Range("E" + $src).Value()
    $ws.Range("C" + $n).Value = $cVal
    $ws.Range("D" + $n).Value = $dVal
    $ws.Range("E" + $n).Value = $eVal
}

# 4) Row 33 becomes the old row 12: restore the stashed formatting then the stashed values.
$ws.Range("A200:E200").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122)
$ws.Range("C33").Value = $origC12
$ws.Range("D33").Value = $origD12
$ws.Range("E33").Value = $origE12

# 5) Drop the scratch row entirely so it doesn't linger in the used range.
$ws.Rows(200).Delete()

Write-Host "Row move complete"
